$d = $word.ActiveDocument

# Step 1: shrink the original run's text down to " method", leaving the
# remainder of the sentence to be rebuilt as a sequence of new runs.
# $findRange is updated in place by Find.Execute to span the replacement
# text (" method"), so its .End is exactly the insertion point we need.
$findRange = $d.Content
$found1 = $findRange.Find.Execute(
    " method. This method writes the final flight path into a new file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " method", 2)

$insertPoint = $findRange.End

$chunks = @(
    " which",
    " ",
    "writes ",
    "the final flight path",
    ", including the number of stops, ",
    "into a new file",
    " ",
    "."
)

$ranges = @()
foreach ($chunk in $chunks) {
    $ins = $d.Range($insertPoint, $insertPoint)
    $ins.Text = $chunk
    $len = $chunk.Length
    $newRun = $d.Range($insertPoint, $insertPoint + $len)
    # Temporarily toggle Bold on so this freshly typed span is kept as its
    # own run instead of being silently re-merged into the previous,
    # identically formatted run.
    $newRun.Bold = 1
    $ranges = $ranges + $newRun
    $insertPoint = $insertPoint + $len
}

# Now turn Bold back off on every chunk we just inserted, one range at a
# time; because each chunk already lives in its own run, clearing Bold
# like this does not cause the runs to coalesce back together.
foreach ($nr in $ranges) {
    $nr.Bold = 0
}
